$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The two old group-name strings ("Nhà 1 tầng, tường 111" / "Nhà 1 tầng, tường 200")
# are both replaced by a single new name used across B2:B6.
$newName = "Giá Thuế Tài Nguyên Khoảng sản kim loại"
$ws.Range("B2").Value = $newName
$ws.Range("B3").Value = $newName
$ws.Range("B4").Value = $newName
$ws.Range("B5").Value = $newName
$ws.Range("B6").Value = $newName

# Widen column B to fit the longer text (resulting stored width = 60 characters).
$ws.Columns("B").ColumnWidth = 59.14

# Move the active cell/selection as recorded in the saved workbook.
$ws.Range("E15").Select() | Out-Null
